# Insert a new data row at row 19 (pushes existing rows 19.. down by one)
# and populate it with the new record's values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows(19).Insert()

$ws.Cells.Item(19, 1).Value = 3
$ws.Cells.Item(19, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(19, 3).Value = "Coquimbo"
$ws.Cells.Item(19, 4).Value = 44547
$ws.Cells.Item(19, 5).Value = 5
$ws.Cells.Item(19, 6).Value = 100112030
$ws.Cells.Item(19, 7).Value = "Poroto granado"
$ws.Cells.Item(19, 8).Value = "Sin especificar"
$ws.Cells.Item(19, 9).Value = "Primera"
$ws.Cells.Item(19, 10).Value = 73
$ws.Cells.Item(19, 11).Value = 41000
$ws.Cells.Item(19, 12).Value = 42000
$ws.Cells.Item(19, 13).Value = 41521
$ws.Cells.Item(19, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(19, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(19, 16).Value = 1661
$ws.Cells.Item(19, 17).Value = 25
$ws.Cells.Item(19, 18).Value = "Hortaliza"

$ws.Cells.Item(19, 4).NumberFormat = $ws.Cells.Item(20, 4).NumberFormat
